# Richtlinien Nachhaltigkeitstagung.docx
# Korrektur die Arbeitspaketen-Name und -Nummerierung.
#
#  1. "Soll im Jahr 2012 Stattfinden" -> "Soll im Jahr 2021 Stattfinden"
#     (only the "12" inside "2012" is touched; the surrounding text is
#     left in place, which is why the saved document ends up with the
#     sentence split across three runs).
#  2. Drop the stray "_GoBack" bookmark that Word leaves behind from the
#     previous editing session.

$d = $word.ActiveDocument

# --- 1. Fix the year -------------------------------------------------
# Find the "2012" in the body text (this only matches the real sentence;
# the only other occurrence of the string "2012" in the package lives in
# an XML namespace URI, which Find never sees).
$yearRange = $d.Content
if ($yearRange.Find.Execute("2012", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {

    # Only replace the last two digits ("12" -> "21"); this mirrors the
    # original edit and keeps "Soll im Jahr 20" intact as its own run.
    $digits = $d.Range($yearRange.Start + 2, $yearRange.End)
    $digits.Text = "21"

    # Force Word to give these two characters their own run instead of
    # silently re-merging them into the neighbouring text (toggling a
    # character property on/off is a no-op visually/semantically, but it
    # is what makes the run boundary "stick").
    $digits.Bold = 1
    $digits.Bold = 0
}

# --- 2. Remove the leftover _GoBack bookmark -------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
